# Reduce to 13 countries: remove BGR, CAN, FIN, LVA, NLD, RUS, SRB rows.
# Original row layout (1-indexed):
#   3 BRA, 4 BGR, 5 CAN, 6 CHL, 7 EST, 8 FIN, 9 GEO, 10 IDN, 11 ITA, 12 LVA,
#   13 LTU, 14 NLD, 15 PER, 16 POL, 17 PRT, 18 RUS, 19 SRB, 20 SVK, 21 ESP, 22 USA
# Delete from the bottom up so earlier row numbers stay valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(19).Delete()
$ws.Rows(18).Delete()
$ws.Rows(14).Delete()
$ws.Rows(12).Delete()
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
$ws.Rows(4).Delete()

# Correct a data value: GEO's "Use" (column K) figure, now on row 6.
$ws.Range("K6").Value = 0.834
